# "add participant to groupe"
# Replace the two e-mail addresses with a new pair, remove the hyperlink
# (and its associated style) that used to sit on A3, and leave the
# selection where the author left it (Q3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values. Order matters for how new strings land in the
# shared-strings table (mirrors the order the original author typed them
# in): A3 first, then A2, then A1 (A1 keeps its existing text "email").
$ws.Range("A3").Value = "hello@yahoo.fr"
$ws.Range("A2").Value = "droitformation.web@gmail.com"
$ws.Range("A1").Value = "email"

# The old A3 value was a mailto: hyperlink with the "Lien hypertexte"
# style; the new participant address is plain text, so drop the
# hyperlink and restore the default cell style.
$ws.Hyperlinks.Delete()
$ws.Range("A3").Style = "Normal"

# Restore the selection/active cell to where it was left in the file.
$ws.Range("Q3").Select() | Out-Null
